$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "26.150.53"
$ws.Range("E2").Value = "  -1.25%  "
$ws.Range("D3").Value = "1.660.86"
$ws.Range("E3").Value = "  -0.79%  "
$ws.Range("E4").Value = "  +0.11%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "216.54"
$ws.Range("E5").Value = "  -1.50%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.5152"
$ws.Range("E6").Value = "  -3.16%  "
$ws.Range("E7").Value = "  +0.10%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.2640"
$ws.Range("E8").Value = "  -1.99%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.06271"
$ws.Range("E9").Value = "  -1.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "20.76"
$ws.Range("E10").Value = "  -4.56%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.07744"
$ws.Range("E11").Value = "  -0.68%  "
$ws.Range("B12").Value = "WrappedEther"
$ws.Range("C12").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D12").Value = "1.659.00"
$ws.Range("E12").Value = "  -0.92%  "
$ws.Range("B13").Value = "Polkadot"
$ws.Range("C13").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "4.435"
$ws.Range("E13").Value = "  -1.41%  "
$ws.Range("D14").Value = "1.886.88"
$ws.Range("E14").Value = "  -0.89%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "0.5426"
$ws.Range("E15").Value = "  -2.60%  "
$ws.Range("D16").Value = "0.0₅8100"
$ws.Range("E16").Value = "  -2.94%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "64.80"
$ws.Range("E17").Value = "  -1.31%  "
$ws.Range("D18").Value = "26.162.56"
$ws.Range("E18").Value = "  -1.28%  "
$ws.Range("E19").Value = "  +0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "4.618"
$ws.Range("E20").Value = "  -2.99%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "191.77"
$ws.Range("E21").Value = "  -0.39%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "10.08"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "6.021"
$ws.Range("E23").Value = "  -4.67%  "
$ws.Range("E24").Value = "  +0.16%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "139.72"
$ws.Range("E25").Value = "  +0.45%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.1221"
$ws.Range("E26").Value = "  -4.75%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.225"
$ws.Range("E27").Value = "  -2.61%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "16.10"
$ws.Range("E28").Value = "  -1.10%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "1.429"
$ws.Range("E29").Value = "  +0.14%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.05952"
$ws.Range("E30").Value = "  -5.16%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.270"
$ws.Range("E31").Value = "  -1.71%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.571"
$ws.Range("E32").Value = "  -0.97%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "3.256"
$ws.Range("E33").Value = "  -4.97%  "
$ws.Range("E34").Value = "  -5.60%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.9636"
$ws.Range("E35").Value = "  -4.68%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "2.420"
$ws.Range("E36").Value = "  -0.10%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.769"
$ws.Range("E37").Value = "  -0.41%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.5651"
$ws.Range("E38").Value = "  -8.11%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01590"
$ws.Range("E39").Value = "  -1.72%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "5.962"
$ws.Range("E40").Value = "  -2.71%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.8555"
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  +0.14%  "
$ws.Range("D43").Value = "1.013.21"
$ws.Range("E43").Value = "  -7.43%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "100.54"
$ws.Range("E44").Value = "  -0.08%  "
$ws.Range("D45").Value = "1.801.16"
$ws.Range("E45").Value = "  -0.99%  "
$ws.Range("E46").Value = "  -2.14%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "56.59"
$ws.Range("E47").Value = "  -3.49%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.006"
$ws.Range("E48").Value = "  +0.31%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.995"
$ws.Range("E49").Value = "  -2.22%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.05165"
$ws.Range("E50").Value = "  -0.57%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "1.452"
$ws.Range("E51").Value = "  -4.19%  "

Write-Host "Applied 97 cell changes"